$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for cryptocurrency rows
# D-column values are forced to remain text (matching the original inlineStr type)
# by temporarily applying a text number format, then restoring the default style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.640.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.514.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.73%  "

$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -2.16%  "

$ws.Range("E10").Value = "  -2.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("E13").Value = "  -2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.900.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.508.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.92%  "

$ws.Range("E17").Value = "  -2.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.654.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.79%  "

$ws.Range("E20").Value = "  -1.71%  "

$ws.Range("E21").Value = "  -0.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("E25").Value = "  -2.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.73%  "

$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("E28").Value = "  +8.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("E33").Value = "  +3.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.14%  "

$ws.Range("E35").Value = "  -3.04%  "

$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("E37").Value = "  -5.29%  "

$ws.Range("E38").Value = "  -1.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.15%  "

$ws.Range("E41").Value = "  +3.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("E43").Value = "  -1.71%  "

$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("E45").Value = "  -1.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.045.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.755.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.90%  "

$ws.Range("E51").Value = "  -0.66%  "
